$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 226, shifting existing rows 226:261 down to 227:262.
$ws.Rows("226:226").Insert()

# Populate the new row 226 with the new weekly price entry.
$ws.Cells.Item(226, 1).Value  = 7
$ws.Cells.Item(226, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(226, 3).Value  = "Ñuble"
$ws.Cells.Item(226, 4).Value  = 45015
$ws.Cells.Item(226, 5).Value  = 16
$ws.Cells.Item(226, 6).Value  = "Fruta"
$ws.Cells.Item(226, 7).Value  = 100104
$ws.Cells.Item(226, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(226, 9).Value  = 100104005
$ws.Cells.Item(226, 10).Value = "Pera"
$ws.Cells.Item(226, 11).Value = "Packham's Triumph"
$ws.Cells.Item(226, 12).Value = "Primera"
$ws.Cells.Item(226, 13).Value = 50
$ws.Cells.Item(226, 14).Value = 10000
$ws.Cells.Item(226, 15).Value = 10000
$ws.Cells.Item(226, 16).Value = 10000
$ws.Cells.Item(226, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(226, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(226, 19).Value = 556
$ws.Cells.Item(226, 20).Value = 18
